# Mirror old PR for User testing to keep some good stuff
# Strip the stray leading space from the "Name.Number" (column E) values
# that were originally generated with a leading blank, and restore the
# selected cell to E2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "Mcbain.358"
$ws.Range("E4").Value = "Huot.633"
$ws.Range("E5").Value = "Bezanson.660"
$ws.Range("E7").Value = "Cavalieri.839"
$ws.Range("E8").Value = "Bass.111"
$ws.Range("E9").Value = "Isaacs.442"
$ws.Range("E10").Value = "Wolanski.136"
$ws.Range("E11").Value = "Clendening.100"
$ws.Range("E12").Value = "Ashbaugh.213"
$ws.Range("E13").Value = "Pinder.964"
$ws.Range("E14").Value = "Depalma.763"
$ws.Range("E15").Value = "Imel.715"
$ws.Range("E16").Value = "Shambaugh.40"
$ws.Range("E17").Value = "Albritton.668"
$ws.Range("E18").Value = "Slocum.139"
$ws.Range("E19").Value = "Maysonet.58"
$ws.Range("E20").Value = "Maag.400"
$ws.Range("E21").Value = "Laguardia.605"
$ws.Range("E22").Value = "Rosso.903"
$ws.Range("E23").Value = "Jeanbaptiste.831"
$ws.Range("E24").Value = "Siegrist.361"
$ws.Range("E25").Value = "Balas.71"
$ws.Range("E26").Value = "Costello.760"
$ws.Range("E27").Value = "Poulter.816"
$ws.Range("E28").Value = "Whitmire.337"
$ws.Range("E29").Value = "Politte.833"
$ws.Range("E30").Value = "Hartshorn.948"
$ws.Range("E31").Value = "Iser.962"
$ws.Range("E32").Value = "Tinajero.70"
$ws.Range("E33").Value = "Caplan.540"
$ws.Range("E34").Value = "Seppala.173"
$ws.Range("E35").Value = "Statler.8"
$ws.Range("E36").Value = "Dartez.850"
$ws.Range("E37").Value = "Backman.359"
$ws.Range("E38").Value = "Norgard.983"
$ws.Range("E39").Value = "Ingwersen.694"
$ws.Range("E41").Value = "Landey.624"
$ws.Range("E42").Value = "Herbst.123"
$ws.Range("E43").Value = "Mings.349"
$ws.Range("E44").Value = "Gandara.979"
$ws.Range("E45").Value = "Caffrey.986"
$ws.Range("E46").Value = "Branning.888"
$ws.Range("E47").Value = "Beland.117"
$ws.Range("E48").Value = "Mariscal.729"
$ws.Range("E49").Value = "Medel.318"
$ws.Range("E50").Value = "Normandin.386"
$ws.Range("E51").Value = "Breuer.491"
$ws.Range("E52").Value = "Liberto.857"
$ws.Range("E53").Value = "Godfrey.949"
$ws.Range("E54").Value = "Mikkelson.973"
$ws.Range("E55").Value = "Fleagle.381"
$ws.Range("E56").Value = "Darlington.481"
$ws.Range("E57").Value = "Hosein.882"
$ws.Range("E58").Value = "Juarbe.11"
$ws.Range("E59").Value = "Deslauriers.767"
$ws.Range("E60").Value = "Horney.105"
$ws.Range("E61").Value = "Nivens.576"
$ws.Range("E62").Value = "Enfinger.107"
$ws.Range("E63").Value = "Mccraw.316"
$ws.Range("E64").Value = "Eggleton.598"
$ws.Range("E65").Value = "Wymore.153"
$ws.Range("E66").Value = "Dayton.843"
$ws.Range("E67").Value = "Petti.427"
$ws.Range("E68").Value = "Guercio.396"
$ws.Range("E69").Value = "Vasconcellos.50"
$ws.Range("E70").Value = "Heal.856"
$ws.Range("E71").Value = "Hornick.940"
$ws.Range("E72").Value = "Franson.807"
$ws.Range("E73").Value = "Mahar.143"
$ws.Range("E74").Value = "Denney.753"
$ws.Range("E75").Value = "Saffell.584"
$ws.Range("E76").Value = "Maze.620"
$ws.Range("E77").Value = "Pawlowicz.398"
$ws.Range("E78").Value = "Lahr.942"
$ws.Range("E79").Value = "Lerman.315"
$ws.Range("E81").Value = "Delagarza.200"
$ws.Range("E82").Value = "Thompkins.432"
$ws.Range("E83").Value = "Mccroskey.8"
$ws.Range("E84").Value = "Wynter.419"
$ws.Range("E85").Value = "Cliff.333"
$ws.Range("E86").Value = "Storck.591"
$ws.Range("E87").Value = "Meikle.556"
$ws.Range("E88").Value = "Lightfoot.559"
$ws.Range("E89").Value = "Ester.921"
$ws.Range("E90").Value = "Mungo.80"
$ws.Range("E91").Value = "Michalski.556"
$ws.Range("E92").Value = "Lunceford.641"
$ws.Range("E93").Value = "Siebert.643"
$ws.Range("E94").Value = "Call.939"
$ws.Range("E95").Value = "Robin.11"
$ws.Range("E96").Value = "Roddy.883"
$ws.Range("E97").Value = "Pangle.362"
$ws.Range("E98").Value = "Hastings.286"
$ws.Range("E99").Value = "Seller.821"
$ws.Range("E100").Value = "Mumaw.240"
$ws.Range("E101").Value = "Popham.241"
$ws.Range("E102").Value = "Grieve.693"
$ws.Range("E103").Value = "Camacho.4"
$ws.Range("E104").Value = "Ehrenberg.441"
$ws.Range("E105").Value = "Champine.202"
$ws.Range("E106").Value = "Wells.55"
$ws.Range("E107").Value = "Spicer.41"
$ws.Range("E108").Value = "Shirley.671"
$ws.Range("E109").Value = "Eddington.306"
$ws.Range("E110").Value = "Morrisette.605"
$ws.Range("E111").Value = "Whiting.851"
$ws.Range("E112").Value = "Nordahl.9"
$ws.Range("E113").Value = "Mullet.321"
$ws.Range("E114").Value = "Rand.840"
$ws.Range("E115").Value = "Mathisen.218"
$ws.Range("E116").Value = "Garten.560"
$ws.Range("E117").Value = "Tyree.611"
$ws.Range("E118").Value = "Peppler.143"
$ws.Range("E119").Value = "Paulson.371"
$ws.Range("E121").Value = "Lickteig.332"
$ws.Range("E122").Value = "Haddon.17"
$ws.Range("E123").Value = "Mauger.736"
$ws.Range("E124").Value = "Berge.171"
$ws.Range("E125").Value = "Willcutt.973"
$ws.Range("E126").Value = "Gaver.489"
$ws.Range("E127").Value = "Ruyle.883"
$ws.Range("E128").Value = "Heeter.44"
$ws.Range("E129").Value = "Fitzwater.145"
$ws.Range("E130").Value = "Aquilino.382"
$ws.Range("E131").Value = "Palomino.241"
$ws.Range("E132").Value = "Dunson.706"
$ws.Range("E133").Value = "Hager.952"
$ws.Range("E134").Value = "Bussell.732"
$ws.Range("E135").Value = "Reese.274"
$ws.Range("E136").Value = "Borjas.930"
$ws.Range("E137").Value = "Cady.595"
$ws.Range("E138").Value = "Okada.299"
$ws.Range("E139").Value = "Debow.528"
$ws.Range("E140").Value = "Calderone.593"
$ws.Range("E141").Value = "Berrier.471"
$ws.Range("E142").Value = "Parry.560"
$ws.Range("E143").Value = "Tignor.161"
$ws.Range("E144").Value = "Cutsforth.174"
$ws.Range("E145").Value = "Cripe.236"
$ws.Range("E146").Value = "Belmont.982"
$ws.Range("E147").Value = "Girard.502"
$ws.Range("E148").Value = "Vining.935"
$ws.Range("E149").Value = "Clever.468"
$ws.Range("E150").Value = "Dias.832"
$ws.Range("E151").Value = "Drake.377"
$ws.Range("E152").Value = "Lane.218"
$ws.Range("E153").Value = "Gerry.184"
$ws.Range("E154").Value = "Younts.552"
$ws.Range("E155").Value = "Mcmurtrie.820"
$ws.Range("E156").Value = "Dowless.286"
$ws.Range("E157").Value = "Scheffer.323"
$ws.Range("E158").Value = "Rempel.454"
$ws.Range("E159").Value = "Lis.287"
$ws.Range("E161").Value = "Lesperance.914"
$ws.Range("E162").Value = "Gunning.855"
$ws.Range("E163").Value = "Cripps.973"
$ws.Range("E164").Value = "Younkin.677"
$ws.Range("E165").Value = "Braz.539"
$ws.Range("E166").Value = "Zielinski.230"
$ws.Range("E167").Value = "Junior.195"
$ws.Range("E168").Value = "Veneziano.661"
$ws.Range("E169").Value = "Down.797"
$ws.Range("E170").Value = "Ogle.891"
$ws.Range("E171").Value = "Resnick.437"
$ws.Range("E172").Value = "Tindell.641"
$ws.Range("E173").Value = "Woltz.563"
$ws.Range("E174").Value = "Borjas.78"
$ws.Range("E175").Value = "Esquer.919"
$ws.Range("E176").Value = "Zdenek.728"
$ws.Range("E177").Value = "Astorga.32"
$ws.Range("E178").Value = "Dollins.210"
$ws.Range("E179").Value = "Standifer.77"
$ws.Range("E180").Value = "Stoke.982"
$ws.Range("E181").Value = "Grief.860"
$ws.Range("E182").Value = "Denis.806"
$ws.Range("E183").Value = "Gamblin.788"
$ws.Range("E184").Value = "Faddis.634"
$ws.Range("E185").Value = "Binns.159"
$ws.Range("E186").Value = "Hepfer.154"
$ws.Range("E187").Value = "Slade.484"
$ws.Range("E188").Value = "Papp.563"
$ws.Range("E189").Value = "Encinas.262"
$ws.Range("E190").Value = "Mcdonald.342"
$ws.Range("E191").Value = "Benigno.116"
$ws.Range("E192").Value = "Brinkley.866"
$ws.Range("E193").Value = "Certain.335"
$ws.Range("E194").Value = "Cypher.353"
$ws.Range("E195").Value = "Lautenschlage.917"
$ws.Range("E196").Value = "Rushford.890"
$ws.Range("E197").Value = "Delossantos.38"
$ws.Range("E198").Value = "Slee.951"
$ws.Range("E199").Value = "Moretti.365"
$ws.Range("E201").Value = "Rainer.108"
$ws.Range("E202").Value = "Breitenstein.677"
$ws.Range("E203").Value = "Coletta.344"
$ws.Range("E204").Value = "Hieber.798"
$ws.Range("E205").Value = "Closson.956"
$ws.Range("E206").Value = "Tully.114"
$ws.Range("E207").Value = "Caywood.170"
$ws.Range("E208").Value = "Silsby.242"
$ws.Range("E209").Value = "Tait.159"
$ws.Range("E210").Value = "Aleman.73"
$ws.Range("E211").Value = "Brindley.439"
$ws.Range("E212").Value = "Dicarlo.674"
$ws.Range("E213").Value = "Florez.613"
$ws.Range("E214").Value = "Sawyers.890"
$ws.Range("E215").Value = "Byrne.461"
$ws.Range("E216").Value = "Hagar.119"
$ws.Range("E217").Value = "Goldie.250"
$ws.Range("E218").Value = "Monfort.301"
$ws.Range("E219").Value = "Toll.679"
$ws.Range("E220").Value = "Gammage.15"
$ws.Range("E221").Value = "Darnell.242"
$ws.Range("E222").Value = "Paro.334"
$ws.Range("E223").Value = "Barwick.128"
$ws.Range("E224").Value = "Dugger.786"
$ws.Range("E225").Value = "Leibowitz.835"
$ws.Range("E226").Value = "Krebsbach.938"

$null = $ws.Range("E2").Select()
